$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ticket 171: quitado ID_SERVICIO de plantillas y puesto ID_ORDEN_SERVICIO
# Replace the ${ID_SERVICIO} placeholder text with ${ID_ORDEN_SERVICIO}
# (single-quoted so PowerShell does not try to expand ${...} as a variable)
$ws.Range("A5").Value2 = '${ID_ORDEN_SERVICIO}'

# The merged A5:A6 header cell also loses its bold weight and gains word
# wrapping (its xf moved from the "fontId=4 / no wrap" entry to the
# "fontId=5 / wrap" entry used by the rest of the row).
$headerCell = $ws.Range("A5:A6")
$headerCell.Font.Bold = $false
$headerCell.WrapText = $true

# Move the active selection from F7 to A5
[void]$ws.Range("A5").Select()
